$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (label "R") updated values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 166
$wsOff.Range("C3").Value = 110
$wsOff.Range("D3").Value = 55
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 4
$wsOff.Range("G3").Value = 2

# Sheet "DEF" - row 3 (label "R") updated values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 176
$wsDef.Range("C3").Value = 122
$wsDef.Range("D3").Value = 38
$wsDef.Range("E3").Value = 17
